# Fruta / hortaliza, semanal
# Insert two new weekly records ("Navel Late" / "Primera" and "Navel Late" /
# "Segunda", dated 45223) right after the existing row 168, pushing the
# remaining historical rows (old 169-187) down to 171-189.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 169-170; everything from the old row 169 onward
# shifts down by two rows (old 169 -> new 171, ..., old 187 -> new 189).
$ws.Range("A169:A170").EntireRow.Insert()

# New row 169: Navel Late, Primera
$ws.Cells.Item(169, 1).Value = 1
$ws.Cells.Item(169, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(169, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(169, 4).Value = 45223
$ws.Cells.Item(169, 5).Value = 15
$ws.Cells.Item(169, 6).Value = "Fruta"
$ws.Cells.Item(169, 7).Value = 100102
$ws.Cells.Item(169, 8).Value = "Cítricos"
$ws.Cells.Item(169, 9).Value = 100102005
$ws.Cells.Item(169, 10).Value = "Naranja"
$ws.Cells.Item(169, 11).Value = "Navel Late"
$ws.Cells.Item(169, 12).Value = "Primera"
$ws.Cells.Item(169, 13).Value = 150
$ws.Cells.Item(169, 14).Value = 750
$ws.Cells.Item(169, 15).Value = 750
$ws.Cells.Item(169, 16).Value = 750
$ws.Cells.Item(169, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(169, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(169, 19).Value = 750
$ws.Cells.Item(169, 20).Value = 1

# New row 170: Navel Late, Segunda
$ws.Cells.Item(170, 1).Value = 1
$ws.Cells.Item(170, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(170, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(170, 4).Value = 45223
$ws.Cells.Item(170, 5).Value = 15
$ws.Cells.Item(170, 6).Value = "Fruta"
$ws.Cells.Item(170, 7).Value = 100102
$ws.Cells.Item(170, 8).Value = "Cítricos"
$ws.Cells.Item(170, 9).Value = 100102005
$ws.Cells.Item(170, 10).Value = "Naranja"
$ws.Cells.Item(170, 11).Value = "Navel Late"
$ws.Cells.Item(170, 12).Value = "Segunda"
$ws.Cells.Item(170, 13).Value = 150
$ws.Cells.Item(170, 14).Value = 800
$ws.Cells.Item(170, 15).Value = 800
$ws.Cells.Item(170, 16).Value = 800
$ws.Cells.Item(170, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(170, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(170, 19).Value = 800
$ws.Cells.Item(170, 20).Value = 1
